# Insert two new data rows (weekly price records) at row 580 of Sheet1,
# pushing the existing rows 580:607 down to 582:609.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank rows before the current row 580, shifting everything
# below (including the old row 580) downward.
$ws.Rows("580:581").Insert()

# Row 580 (new): 2023-08-09 record for "Región Metropolitana"
$ws.Range("A580").Value = 5
$ws.Range("B580").Value = "Macroferia Regional de Talca"
$ws.Range("C580").Value = "Maule"
$ws.Range("D580").Value = 45147
$ws.Range("D580").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E580").Value = 7
$ws.Range("F580").Value = 100112023
$ws.Range("G580").Value = "Brócoli"
$ws.Range("H580").Value = "Sin especificar"
$ws.Range("I580").Value = "Primera"
$ws.Range("J580").Value = 3000
$ws.Range("K580").Value = 900
$ws.Range("L580").Value = 900
$ws.Range("M580").Value = 900
$ws.Range("N580").Value = "$/unidad"
$ws.Range("O580").Value = "Región Metropolitana"
$ws.Range("P580").Value = 900
$ws.Range("Q580").Value = 1
$ws.Range("R580").Value = "Hortaliza"

# Row 581 (new): same date, but "Región del Maule"
$ws.Range("A581").Value = 5
$ws.Range("B581").Value = "Macroferia Regional de Talca"
$ws.Range("C581").Value = "Maule"
$ws.Range("D581").Value = 45147
$ws.Range("D581").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E581").Value = 7
$ws.Range("F581").Value = 100112023
$ws.Range("G581").Value = "Brócoli"
$ws.Range("H581").Value = "Sin especificar"
$ws.Range("I581").Value = "Primera"
$ws.Range("J581").Value = 3000
$ws.Range("K581").Value = 800
$ws.Range("L581").Value = 800
$ws.Range("M581").Value = 800
$ws.Range("N581").Value = "$/unidad"
$ws.Range("O581").Value = "Región del Maule"
$ws.Range("P581").Value = 800
$ws.Range("Q581").Value = 1
$ws.Range("R581").Value = "Hortaliza"
